$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain stored as text so values such as
# "1.00", "64.844.38" or "0.778" are not re-interpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.844.38"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.457.65"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "576.25"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "160.90"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  +10.93%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "3.463.05"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +3.87%  "
$ws.Range("D13").Value = "4.047.43"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "28.36"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "64.827.75"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "3.471.62"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "381.57"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "0.553"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").Value = "72.78"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +9.40%  "
$ws.Range("D31").Value = "6.21"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").Value = "23.57"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "7.27"
$ws.Range("E34").Value = "  +6.28%  "
$ws.Range("D35").Value = "1.62"
$ws.Range("E35").Value = "  +11.21%  "
$ws.Range("D36").Value = "160.84"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  +5.17%  "
$ws.Range("D38").Value = "0.0779"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "2.938.04"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "6.77"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.68"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "4.66"
$ws.Range("E42").Value = "  +7.29%  "
$ws.Range("D43").Value = "0.0322"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "42.81"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.778"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "26.01"
$ws.Range("E46").Value = "  +11.77%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.111"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "317.53"
$ws.Range("E49").Value = "  +8.65%  "
$ws.Range("D50").Value = "0.879"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "6.62"
$ws.Range("E51").Value = "  +3.68%  "
